# Update odds values in row 4 (hCptA7hl - Correcaminos vs Atl. Morelia)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.67
$ws.Range("I4").Value = 2.55
$ws.Range("J4").Value = 3.3
$ws.Range("L4").Value = 3.15
$ws.Range("N4").Value = 1.06
$ws.Range("P4").Value = 2.95
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.78
$ws.Range("T4").Value = 2.47
$ws.Range("W4").Value = 8.5
$ws.Range("Y4").Value = 9.75
$ws.Range("Z4").Value = 32
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 30
$ws.Range("AH4").Value = 8.75
$ws.Range("AJ4").Value = 9.5
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 28
$ws.Range("AN4").Value = 4.55
$ws.Range("AO4").Value = 14.5
$ws.Range("AP4").Value = 22
$ws.Range("AQ4").Value = 65
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 300
$ws.Range("AT4").Value = 2.45
$ws.Range("AU4").Value = 6.7
$ws.Range("AW4").Value = 4.45
$ws.Range("AX4").Value = 14
$ws.Range("AY4").Value = 21
$ws.Range("AZ4").Value = 60
$ws.Range("BA4").Value = 90
$ws.Range("BB4").Value = 250
